$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '65.875.44'
Set-TextCell 2 5 '  +0.80%  '

Set-TextCell 3 4 '3.387.55'
Set-TextCell 3 5 '  -0.54%  '

Set-TextCell 4 4 '0.999'
Set-TextCell 4 5 '  -0.08%  '

Set-TextCell 5 4 '564.28'
Set-TextCell 5 5 '  +0.54%  '

Set-TextCell 6 4 '176.54'
Set-TextCell 6 5 '  +0.32%  '

Set-TextCell 7 5 '  +0.48%  '

Set-TextCell 8 4 '3.380.66'
Set-TextCell 8 5 '  -0.42%  '

Set-TextCell 9 5 '  -0.06%  '

Set-TextCell 10 5 '  +1.94%  '

Set-TextCell 11 4 '0.634'
Set-TextCell 11 5 '  +0.30%  '

Set-TextCell 12 4 '53.87'
Set-TextCell 12 5 '  -1.97%  '

Set-TextCell 13 4 '0.0000279'
Set-TextCell 13 5 '  -0.62%  '

Set-TextCell 14 4 '9.24'
Set-TextCell 14 5 '  +0.81%  '

Set-TextCell 15 4 '3.925.25'
Set-TextCell 15 5 '  -0.63%  '

Set-TextCell 16 2 'Chainlink'
Set-TextCell 16 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 16 4 '18.23'
Set-TextCell 16 5 '  -0.97%  '

Set-TextCell 17 2 'TRON'
Set-TextCell 17 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 17 4 '0.120'
Set-TextCell 17 5 '  +0.35%  '

Set-TextCell 18 4 '3.380.33'
Set-TextCell 18 5 '  -0.59%  '

Set-TextCell 19 4 '65.885.23'
Set-TextCell 19 5 '  +0.77%  '

Set-TextCell 20 4 '11.89'
Set-TextCell 20 5 '  -0.08%  '

Set-TextCell 21 4 '0.996'
Set-TextCell 21 5 '  -0.10%  '

Set-TextCell 22 4 '464.82'
Set-TextCell 22 5 '  -1.66%  '

Set-TextCell 23 4 '4.93'
Set-TextCell 23 5 '  -1.78%  '

Set-TextCell 24 4 '14.53'
Set-TextCell 24 5 '  +7.75%  '

Set-TextCell 25 4 '89.66'
Set-TextCell 25 5 '  +2.65%  '

Set-TextCell 26 5 '  -1.15%  '

Set-TextCell 27 5 '  +0.12%  '

Set-TextCell 28 4 '10.64'
Set-TextCell 28 5 '  -2.76%  '

Set-TextCell 29 4 '8.71'
Set-TextCell 29 5 '  -1.62%  '

Set-TextCell 30 4 '31.14'
Set-TextCell 30 5 '  -0.66%  '

Set-TextCell 31 4 '6.60'
Set-TextCell 31 5 '  -2.26%  '

Set-TextCell 32 4 '11.47'
Set-TextCell 32 5 '  -0.90%  '

Set-TextCell 33 4 '581.14'
Set-TextCell 33 5 '  +0.89%  '

Set-TextCell 34 4 '62.30'
Set-TextCell 34 5 '  +0.30%  '

Set-TextCell 35 5 '  -0.39%  '

Set-TextCell 36 5 '  +0.10%  '

Set-TextCell 37 4 '3.60'
Set-TextCell 37 5 '  +1.36%  '

Set-TextCell 38 5 '  +1.31%  '

Set-TextCell 39 4 '36.02'
Set-TextCell 39 5 '  +0.19%  '

Set-TextCell 40 5 '  +1.17%  '

Set-TextCell 41 4 '0.0₃0746'
Set-TextCell 41 5 '  -2.75%  '

Set-TextCell 42 4 '3.104.57'
Set-TextCell 42 5 '  +0.06%  '

Set-TextCell 43 4 '2.85'
Set-TextCell 43 5 '  -0.87%  '

Set-TextCell 44 4 '0.0418'
Set-TextCell 44 5 '  -0.11%  '

Set-TextCell 45 5 '  -1.07%  '

Set-TextCell 46 5 '  -1.65%  '

Set-TextCell 47 4 '3.17'
Set-TextCell 47 5 '  +0.00%  '

Set-TextCell 48 4 '0.999'
Set-TextCell 48 5 '  -0.04%  '

Set-TextCell 49 4 '140.90'
Set-TextCell 49 5 '  +2.41%  '

Set-TextCell 50 4 '8.51'
Set-TextCell 50 5 '  +2.78%  '

Set-TextCell 51 2 'LidoDAOToken'
Set-TextCell 51 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 51 4 '3.16'
Set-TextCell 51 5 '  +9.12%  '
